$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.097.22'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.57%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.579.12'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.80%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.20%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '207.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +11.79%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '564.63'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.96%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.578.95'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.85%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.612'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("E9").Value = '  -0.25%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.678'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.39%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '60.71'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +10.22%  '

$ws.Range("E12").Value = '  -1.23%  '

$ws.Range("E13").Value = '  +7.34%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.23'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.75%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.155.56'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.82%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.586.03'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.58%  '

$ws.Range("E17").Value = '  +0.88%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.86'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.73%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '67.852.36'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.22%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.17'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.36%  '

$ws.Range("E21").Value = '  +0.33%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '401.57'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.67%  '

$ws.Range("B23").Value = 'RenderToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.49'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +12.40%  '

$ws.Range("B24").Value = 'PancakeSwap'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.15'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.33%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.22'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.04%  '

$ws.Range("E26").Value = '  -1.11%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.40'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.51%  '

$ws.Range("E28").Value = '  +9.20%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.16'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.28%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.68'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.50%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '31.38'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.11%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '661.92'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.41%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '12.04'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.41%  '

$ws.Range("B34").Value = 'OKB'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '63.24'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.16%  '

$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.112'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.44%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '40.95'
$ws.Range("D36").Style = "Normal"

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.408'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.34%  '

$ws.Range("E38").Value = '  +0.11%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.27'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +12.45%  '

$ws.Range("E40").Value = '  +0.87%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.166.60'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.15%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.132'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.30%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.08%  '

$ws.Range("E44").Value = '  +2.40%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.79'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +11.81%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0409'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.03%  '

$ws.Range("E47").Value = '  +0.52%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.64'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.21%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.59'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +10.90%  '

$ws.Range("E50").Value = '  -0.44%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '138.20'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.57%  '
